$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Slide 7 ("pip3 install ... robotframework-selenium2screenshots"):
#    remove the duplicated "install " run and add a new paragraph
#    "pip3 install SeleniumLibrary".
# ---------------------------------------------------------------------------
$installSlide = $p.Slides.Item(7)
$contentTr = $installSlide.Shapes.Item(2).TextFrame.TextRange

$fullText = $contentTr.Text
$marker = "install install"
$idx = $fullText.IndexOf($marker)
if ($idx -ge 0) {
    # Delete the first "install " occurrence (characters are 1-based in the COM API).
    $deleteStart = $idx + 1
    $deleteLen = 8  # "install "
    $contentTr.Characters($deleteStart, $deleteLen).Text = ""
}

$contentTr.InsertAfter("`r")
$afterBreak = $installSlide.Shapes.Item(2).TextFrame.TextRange
$afterBreak.InsertAfter("pip3 install ")
$afterBreak2 = $installSlide.Shapes.Item(2).TextFrame.TextRange
$afterBreak2.InsertAfter("SeleniumLibrary")

# ---------------------------------------------------------------------------
# 2. Add two new "Run script with tags" slides right before the final
#    "robot tests..." slide, and push that original slide to the end.
# ---------------------------------------------------------------------------
$lastSlide = $p.Slides.Item($p.Slides.Count)
$lastIndex = $lastSlide.SlideIndex

$dup1 = $lastSlide.Duplicate().Item(1)   # lands right after $lastSlide
$dup2 = $lastSlide.Duplicate().Item(1)   # lands right after $lastSlide, before $dup1

$lastSlide.MoveTo($p.Slides.Count)

$browserTagsSlide = $p.Slides.Item($lastIndex)
$titleTr = $browserTagsSlide.Shapes.Item(1).TextFrame.TextRange
$titleTr.Text = "Run script with tags"
$bodyTr = $browserTagsSlide.Shapes.Item(2).TextFrame.TextRange
$bodyTr.Text = "pybot --variable BROWSER:IE facebook.robot`r"

$functionalTestSlide = $p.Slides.Item($lastIndex + 1)
$titleTr2 = $functionalTestSlide.Shapes.Item(1).TextFrame.TextRange
$titleTr2.Text = "Run script with tags"
$bodyTr2 = $functionalTestSlide.Shapes.Item(2).TextFrame.TextRange
$bodyTr2.Text = "pybot -i 'Functional Test' facebook.robot`r"
